$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 42 new rows starting at row 203 (the Rusturf Tunnel grunts
# Pokemon row plus all of the new Route 116 trainer blocks), matching the
# "insertRow at 203, count 42" structural edit.
$ws.Rows("203:244").Insert()

# Newly inserted rows inherit the formatting of the row above (202);
# normalize them to the sheets plain/default cell style by copying the
# format from a known plain cell (A25) across the whole new block.
$ws.Range("A25:E25").Copy()
$ws.Range("A203:E244").PasteSpecial(-4122)

# Row 202 used to hold the final "END" marker. Repurpose it as the header
# of the new Rusturf Tunnel grunt entry, dropping the special end-of-table
# styling to match the rest of the table.
$ws.Range("A25").Copy()
$ws.Range("A202").PasteSpecial(-4122)
$ws.Range("A202").Value = "TRAINER_GRUNT_RUSTURF_TUNNEL"

$ws.Range("A203").Value = "species"
$ws.Range("B203").Value = "lvl"
$ws.Range("C203").Value = "iv"
$ws.Range("D203").Value = "heldItem"
$ws.Range("E203").Value = "moves"

$ws.Range("A204").Value = "Charmeleon"
$ws.Range("B204").Value = 18

$ws.Range("A206").Value = "# Route 116"

$ws.Range("A208").Value = "TRAINER_JOEY"

$ws.Range("A209").Value = "species"
$ws.Range("B209").Value = "lvl"
$ws.Range("C209").Value = "iv"
$ws.Range("D209").Value = "heldItem"
$ws.Range("E209").Value = "moves"

$ws.Range("A210").Value = "Aron"
$ws.Range("B210").Value = 15

$ws.Range("A211").Value = "Nidoran"
$ws.Range("B211").Value = 16

$ws.Range("A213").Value = "TRAINER_JOSE"

$ws.Range("A214").Value = "species"
$ws.Range("B214").Value = "lvl"
$ws.Range("C214").Value = "iv"
$ws.Range("D214").Value = "heldItem"
$ws.Range("E214").Value = "moves"

$ws.Range("A215").Value = "Ledian"
$ws.Range("B215").Value = 18

$ws.Range("A217").Value = "TRAINER_KAREN_1"

$ws.Range("A218").Value = "species"
$ws.Range("B218").Value = "lvl"
$ws.Range("C218").Value = "iv"
$ws.Range("D218").Value = "heldItem"
$ws.Range("E218").Value = "moves"

$ws.Range("A219").Value = "Nidoran"
$ws.Range("B219").Value = 16

$ws.Range("A220").Value = "Eevee"
$ws.Range("B220").Value = 16

$ws.Range("A222").Value = "TRAINER_CLARK"

$ws.Range("A223").Value = "species"
$ws.Range("B223").Value = "lvl"
$ws.Range("C223").Value = "iv"
$ws.Range("D223").Value = "heldItem"
$ws.Range("E223").Value = "moves"

$ws.Range("A224").Value = "Rhyhorn"
$ws.Range("B224").Value = 16

$ws.Range("A225").Value = "Mawile"
$ws.Range("B225").Value = 17

$ws.Range("A227").Value = "TRAINER_JOHNSON"

$ws.Range("A228").Value = "species"
$ws.Range("B228").Value = "lvl"
$ws.Range("C228").Value = "iv"
$ws.Range("D228").Value = "heldItem"
$ws.Range("E228").Value = "moves"

$ws.Range("A229").Value = "Grimey"
$ws.Range("B229").Value = 17

$ws.Range("A230").Value = "Koffing"
$ws.Range("B230").Value = 17

$ws.Range("A232").Value = "TRAINER_DEVAN"

$ws.Range("A233").Value = "species"
$ws.Range("B233").Value = "lvl"
$ws.Range("C233").Value = "iv"
$ws.Range("D233").Value = "heldItem"
$ws.Range("E233").Value = "moves"

$ws.Range("A234").Value = "Alolan_Geodude"
$ws.Range("B234").Value = 16

$ws.Range("A235").Value = "Larvitar"
$ws.Range("B235").Value = 17

$ws.Range("A237").Value = "TRAINER_SARAH"

$ws.Range("A238").Value = "species"
$ws.Range("B238").Value = "lvl"
$ws.Range("C238").Value = "iv"
$ws.Range("D238").Value = "heldItem"
$ws.Range("E238").Value = "moves"

$ws.Range("A239").Value = "Galarian_Minicorn"
$ws.Range("B239").Value = 19

$ws.Range("A240").Value = "Alolan_Meowth"
$ws.Range("B240").Value = 20

$ws.Range("A242").Value = "TRAINER_DAWSON"

$ws.Range("A243").Value = "species"
$ws.Range("B243").Value = "lvl"
$ws.Range("C243").Value = "iv"
$ws.Range("D243").Value = "heldItem"
$ws.Range("E243").Value = "moves"

$ws.Range("A244").Value = "Delcatty"
$ws.Range("B244").Value = 21

$ws.Range("A246").Value = "TRAINER_JERRY_1"

$ws.Range("A247").Value = "species"
$ws.Range("B247").Value = "lvl"
$ws.Range("C247").Value = "iv"
$ws.Range("D247").Value = "heldItem"
$ws.Range("E247").Value = "moves"

$ws.Range("A248").Value = "Slakoth"
$ws.Range("B248").Value = 23

$ws.Range("A249").Value = "Slowpoke"
$ws.Range("B249").Value = 21

$ws.Range("A251").Value = "TRAINER_JANICE"

$ws.Range("A252").Value = "species"
$ws.Range("B252").Value = "lvl"
$ws.Range("C252").Value = "iv"
$ws.Range("D252").Value = "heldItem"
$ws.Range("E252").Value = "moves"

$ws.Range("A253").Value = "Torchic"
$ws.Range("B253").Value = 22

$ws.Range("A254").Value = "Pikachu"
$ws.Range("B254").Value = 23

# Row 265 is the new final "END" sentinel row; give it the same
# end-of-table styling that row 202 used to carry (copied from A1).
$ws.Range("A1").Copy()
$ws.Range("A265").PasteSpecial(-4122)
$ws.Range("A265").Value = "END"

Write-Host "Applied Route 116 / Rusturf Tunnel grunt trainer data"